$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4819.609
$ws.Range("I40").Value = 4760.75
$ws.Range("J40").Value = 4851
$ws.Range("K40").Value = 4760.75
$ws.Range("L40").Value = 4851
$ws.Range("M40").Value = -4585.75
$ws.Range("N40").Value = -5201

$ws.Range("H70").Value = 1001
$ws.Range("I70").Value = 502
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 1506
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -1236
$ws.Range("N70").Value = -5040

$ws.Range("H73").Value = 1001
$ws.Range("I73").Value = 502
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 1506
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -570
$ws.Range("N73").Value = -6372

$ws.Range("H118").Value = 1080.5714
$ws.Range("I118").Value = 1216
$ws.Range("J118").Value = 900
$ws.Range("K118").Value = 3648
$ws.Range("L118").Value = 2700
$ws.Range("M118").Value = -1991
$ws.Range("N118").Value = -6014

$ws.Range("H135").Value = 629.3913
$ws.Range("I135").Value = 277.8125
$ws.Range("J135").Value = 1433
$ws.Range("K135").Value = 2500.3125
$ws.Range("L135").Value = 12897
$ws.Range("M135").Value = 34.6875
$ws.Range("N135").Value = -17967

$ws.Range("H137").Value = 6501.85
$ws.Range("I137").Value = 1557.3334
$ws.Range("J137").Value = 51002.5
$ws.Range("K137").Value = 4672.0002
$ws.Range("L137").Value = 153007.5
$ws.Range("M137").Value = -2122.0002
$ws.Range("N137").Value = -158107.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 580.4483
$ws.Range("I2").Value = 549.7083
$ws.Range("J2").Value = 728
$ws.Range("K2").Value = 549.7083
$ws.Range("L2").Value = 728
$ws.Range("M2").Value = -436.7083
$ws.Range("N2").Value = -954

$ws.Range("H21").Value = 2833
$ws.Range("I21").Value = 1999
$ws.Range("J21").Value = 3250
$ws.Range("K21").Value = 1999
$ws.Range("L21").Value = 3250
$ws.Range("M21").Value = -1625
$ws.Range("N21").Value = -3998

$ws.Range("H30").Value = 493.75
$ws.Range("I30").Value = 490
$ws.Range("J30").Value = 495
$ws.Range("K30").Value = 490
$ws.Range("L30").Value = 495
$ws.Range("M30").Value = -340
$ws.Range("N30").Value = -795

$ws.Range("H32").Value = 6232.345
$ws.Range("I32").Value = 4839.0933
$ws.Range("J32").Value = 11224.833
$ws.Range("K32").Value = 4839.0933
$ws.Range("L32").Value = 11224.833
$ws.Range("M32").Value = -4552.0933
$ws.Range("N32").Value = -11798.833

$ws.Range("H61").Value = 8936.23
$ws.Range("I61").Value = 2814
$ws.Range("J61").Value = 11657.223
$ws.Range("K61").Value = 2814
$ws.Range("L61").Value = 11657.223
$ws.Range("M61").Value = -2602
$ws.Range("N61").Value = -12081.223

$ws.Range("H74").Value = 163976.66
$ws.Range("I74").Value = 279722.3
$ws.Range("J74").Value = 9649.134
$ws.Range("K74").Value = 279722.3
$ws.Range("L74").Value = 9649.134
$ws.Range("M74").Value = -278848.3
$ws.Range("N74").Value = -11397.134

$ws.Range("H77").Value = 163976.66
$ws.Range("I77").Value = 279722.3
$ws.Range("J77").Value = 9649.134
$ws.Range("K77").Value = 1398611.5
$ws.Range("L77").Value = 48245.67
$ws.Range("M77").Value = -1394243.5
$ws.Range("N77").Value = -56981.67

$ws.Range("H110").Value = 7017.3335
$ws.Range("I110").Value = 8892.200000000001
$ws.Range("J110").Value = 5678.143
$ws.Range("K110").Value = 8892.200000000001
$ws.Range("L110").Value = 5678.143
$ws.Range("M110").Value = -6847.200000000001
$ws.Range("N110").Value = -9768.143

$ws.Range("H116").Value = 580.4483
$ws.Range("I116").Value = 549.7083
$ws.Range("J116").Value = 728
$ws.Range("K116").Value = 549.7083
$ws.Range("L116").Value = 728
$ws.Range("M116").Value = 1744.2917
$ws.Range("N116").Value = -5316

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H136").Value = 8936.23
$ws.Range("I136").Value = 2814
$ws.Range("J136").Value = 11657.223
$ws.Range("K136").Value = 8442
$ws.Range("L136").Value = 34971.669
$ws.Range("M136").Value = -5892
$ws.Range("N136").Value = -40071.669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 580.4483
$ws.Range("I3").Value = 549.7083
$ws.Range("J3").Value = 728
$ws.Range("K3").Value = 549.7083
$ws.Range("L3").Value = 728
$ws.Range("M3").Value = -435.7083
$ws.Range("N3").Value = -956

$ws.Range("H86").Value = 4459.579
$ws.Range("I86").Value = 4173.7856
$ws.Range("J86").Value = 5259.8
$ws.Range("K86").Value = 4173.7856
$ws.Range("L86").Value = 5259.8
$ws.Range("M86").Value = -3050.7856
$ws.Range("N86").Value = -7505.8

$ws.Range("H89").Value = 4459.579
$ws.Range("I89").Value = 4173.7856
$ws.Range("J89").Value = 5259.8
$ws.Range("K89").Value = 20868.928
$ws.Range("L89").Value = 26299
$ws.Range("M89").Value = -15252.928
$ws.Range("N89").Value = -37531

$ws.Range("H134").Value = 1869.659
$ws.Range("I134").Value = 1506.579
$ws.Range("J134").Value = 4169.1665
$ws.Range("K134").Value = 4519.737
$ws.Range("L134").Value = 12507.4995
$ws.Range("M134").Value = -1984.737
$ws.Range("N134").Value = -17577.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4129.979
$ws.Range("I31").Value = 3602.3914
$ws.Range("J31").Value = 4615.36
$ws.Range("K31").Value = 3602.3914
$ws.Range("L31").Value = 4615.36
$ws.Range("M31").Value = -3307.3914
$ws.Range("N31").Value = -5205.36

$ws.Range("H34").Value = 4129.979
$ws.Range("I34").Value = 3602.3914
$ws.Range("J34").Value = 4615.36
$ws.Range("K34").Value = 3602.3914
$ws.Range("L34").Value = 4615.36
$ws.Range("M34").Value = -3400.3914
$ws.Range("N34").Value = -5019.36

$ws.Range("H58").Value = 4261.4116
$ws.Range("I58").Value = 4099.6
$ws.Range("J58").Value = 4328.8335
$ws.Range("K58").Value = 4099.6
$ws.Range("L58").Value = 4328.8335
$ws.Range("M58").Value = -3896.6
$ws.Range("N58").Value = -4734.8335

$ws.Range("H92").Value = 42994.168
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 42994.168
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 42994.168
$ws.Range("N92").Value = -47986.168

$ws.Range("H122").Value = 4356.857
$ws.Range("I122").Value = 3958.2727
$ws.Range("J122").Value = 5818.3335
$ws.Range("K122").Value = 11874.8181
$ws.Range("L122").Value = 17455.0005
$ws.Range("M122").Value = -9424.8181
$ws.Range("N122").Value = -22355.0005

$ws.Range("H134").Value = 3043.8262
$ws.Range("I134").Value = 2542.2942
$ws.Range("J134").Value = 4464.8335
$ws.Range("K134").Value = 7626.882599999999
$ws.Range("L134").Value = 13394.5005
$ws.Range("M134").Value = -5091.882599999999
$ws.Range("N134").Value = -18464.5005

$ws.Range("H136").Value = 4261.4116
$ws.Range("I136").Value = 4099.6
$ws.Range("J136").Value = 4328.8335
$ws.Range("K136").Value = 12298.8
$ws.Range("L136").Value = 12986.5005
$ws.Range("M136").Value = -9748.800000000001
$ws.Range("N136").Value = -18086.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3044.08
$ws.Range("I5").Value = 635.3333
$ws.Range("J5").Value = 3804.7368
$ws.Range("K5").Value = 1905.9999
$ws.Range("L5").Value = 11414.2104
$ws.Range("M5").Value = -1793.9999
$ws.Range("N5").Value = -11638.2104

$ws.Range("H55").Value = 8213
$ws.Range("I55").Value = 900
$ws.Range("J55").Value = 9257.714
$ws.Range("K55").Value = 2700
$ws.Range("L55").Value = 27773.142
$ws.Range("M55").Value = -2523
$ws.Range("N55").Value = -28127.142

$ws.Range("H64").Value = 3872.2222
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3872.2222
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 11616.6666
$ws.Range("N64").Value = -12156.6666

$ws.Range("H67").Value = 3872.2222
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3872.2222
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 11616.6666
$ws.Range("N67").Value = -13488.6666

$ws.Range("H97").Value = 1000454.2
$ws.Range("I97").Value = 5000000
$ws.Range("J97").Value = 567.75
$ws.Range("K97").Value = 15000000
$ws.Range("L97").Value = 1703.25
$ws.Range("M97").Value = -14999504
$ws.Range("N97").Value = -2695.25

$ws.Range("H135").Value = 3044.08
$ws.Range("I135").Value = 635.3333
$ws.Range("J135").Value = 3804.7368
$ws.Range("K135").Value = 5717.9997
$ws.Range("L135").Value = 34242.6312
$ws.Range("M135").Value = -3182.9997
$ws.Range("N135").Value = -39312.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 584.8
$ws.Range("I2").Value = 652.4167
$ws.Range("J2").Value = 314.33334
$ws.Range("K2").Value = 652.4167
$ws.Range("L2").Value = 314.33334
$ws.Range("M2").Value = -539.4167
$ws.Range("N2").Value = -540.33334

$ws.Range("H102").Value = 1236.6666
$ws.Range("I102").Value = 979.9655
$ws.Range("J102").Value = 2300.1428
$ws.Range("K102").Value = 979.9655
$ws.Range("L102").Value = 2300.1428
$ws.Range("M102").Value = 642.0345
$ws.Range("N102").Value = -5544.1428

$ws.Range("H132").Value = 2059.6155
$ws.Range("I132").Value = 1585.5483
$ws.Range("J132").Value = 3896.625
$ws.Range("K132").Value = 4756.644899999999
$ws.Range("L132").Value = 11689.875
$ws.Range("M132").Value = -2226.644899999999
$ws.Range("N132").Value = -16749.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 56492
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 56492
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 56492
$ws.Range("N54").Value = -57780

$ws.Range("H136").Value = 4914.3794
$ws.Range("I136").Value = 4995.3125
$ws.Range("J136").Value = 4814.769
$ws.Range("K136").Value = 14985.9375
$ws.Range("L136").Value = 14444.307
$ws.Range("M136").Value = -12435.9375
$ws.Range("N136").Value = -19544.307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 250000000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 250000000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 750000000
$ws.Range("N122").Value = -750004900
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 11908479
$ws.Range("I132").Value = 15877170
$ws.Range("J132").Value = 2407
$ws.Range("K132").Value = 47631510
$ws.Range("L132").Value = 7221
$ws.Range("M132").Value = -47628980
$ws.Range("N132").Value = -12281
